$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.729.60"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.543.44"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "598.58"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").Value = "135.43"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "3.540.51"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "6.90"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").Value = "4.137.45"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "3.547.14"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "26.95"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "64.633.55"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "10.02"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("E20").Value = "  +5.48%  "
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("D22").Value = "386.07"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +5.84%  "
$ws.Range("D24").Value = "3.680.88"
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("D25").Value = "74.15"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +10.11%  "
$ws.Range("E28").Value = "  +6.08%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("D32").Value = "3.545.84"
$ws.Range("E32").Value = "  +3.00%  "
$ws.Range("E33").Value = "  +20.93%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "23.94"
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "0.143"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").Value = "169.68"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("E39").Value = "  +5.32%  "
$ws.Range("D40").Value = "4.99"
$ws.Range("E40").Value = "  +6.93%  "
$ws.Range("D41").Value = "0.0802"
$ws.Range("E41").Value = "  +5.65%  "
$ws.Range("D42").Value = "0.824"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("E43").Value = "  +17.61%  "
$ws.Range("D44").Value = "42.56"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  +3.60%  "
$ws.Range("E47").Value = "  +8.71%  "
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("E49").Value = "  +5.98%  "
$ws.Range("D50").Value = "2.447.40"
$ws.Range("E50").Value = "  +11.44%  "
$ws.Range("E51").Value = "  +13.76%  "
